$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert first new weekly record before current row 134.
# Everything from row 134 downward shifts down by one row.
$ws.Rows.Item(134).Insert()
$ws.Range("A134").Value = 8
$ws.Range("B134").Value = "Terminal La Palmera de La Serena"
$ws.Range("C134").Value = "Coquimbo"
$ws.Range("D134").Value = 44818
$ws.Range("E134").Value = 4
$ws.Range("F134").Value = 100112037
$ws.Range("G134").Value = "Cebollín"
$ws.Range("H134").Value = "Sin especificar"
$ws.Range("I134").Value = "Primera"
$ws.Range("J134").Value = 1200
$ws.Range("K134").Value = 1400
$ws.Range("L134").Value = 1600
$ws.Range("M134").Value = 1500
$ws.Range("N134").Value = "`$/paquete 6 unidades"
$ws.Range("O134").Value = "Provincia del Elquí"
$ws.Range("P134").Value = 250
$ws.Range("Q134").Value = 6
$ws.Range("R134").Value = "Hortaliza"

# Insert second new weekly record before (post-shift) row 203.
# Everything from row 203 downward shifts down by one row again.
$ws.Rows.Item(203).Insert()
$ws.Range("A203").Value = 8
$ws.Range("B203").Value = "Terminal La Palmera de La Serena"
$ws.Range("C203").Value = "Coquimbo"
$ws.Range("D203").Value = 44816
$ws.Range("E203").Value = 4
$ws.Range("F203").Value = 100112037
$ws.Range("G203").Value = "Cebollín"
$ws.Range("H203").Value = "Sin especificar"
$ws.Range("I203").Value = "Primera"
$ws.Range("J203").Value = 1100
$ws.Range("K203").Value = 1400
$ws.Range("L203").Value = 1600
$ws.Range("M203").Value = 1500
$ws.Range("N203").Value = "`$/paquete 6 unidades"
$ws.Range("O203").Value = "Provincia del Elquí"
$ws.Range("P203").Value = 250
$ws.Range("Q203").Value = 6
$ws.Range("R203").Value = "Hortaliza"
